# Update "Timetsheet - Yingyao Lu.xlsx":
#  - Week7 / Week8: tidy up selection + fix the "Week #" header cell (H3)
#  - Week9: fill in the week's time-log entries
#  - Add two new sheets Week10 and Week11 (cloned from the Week9 template)
#    with their own time-log entries, Print_Area and Week_Start defined
#    names, and make Week10 the active/selected tab.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Week7: only the selection + header-number needed a touch-up
# ---------------------------------------------------------------------------
$ws7 = $wb.Worksheets.Item("Week7")
$ws7.Range("H3").Value = 7
$ws7.Range("H4").Select()

# ---------------------------------------------------------------------------
# Week8: same kind of touch-up (and it loses the "tabSelected" flag once
# Week10 becomes the active tab later on)
# ---------------------------------------------------------------------------
$ws8 = $wb.Worksheets.Item("Week8")
$ws8.Range("H3").Value = 8
$ws8.Range("H4").Select()

# ---------------------------------------------------------------------------
# Week9: populate the timesheet rows for that week
# ---------------------------------------------------------------------------
$ws9 = $wb.Worksheets.Item("Week9")
$ws9.Range("H3").Value = 9

$ws9.Range("C6").Value = 0.41666666666666669
$ws9.Range("D6").Value = 0.625
$ws9.Range("E6").Value = 5
$ws9.Range("F6").Value = "Work on database"
$ws9.Range("G6").Value = "milestone 2 task"

$ws9.Range("C7").Value = 0.54166666666666663
$ws9.Range("D7").Value = 0.875
$ws9.Range("E7").Value = 7
$ws9.Range("F7").Value = "Develop Timeline page"
$ws9.Range("G7").Value = "milestone 2 task"

$ws9.Range("C8").Value = 0.375
$ws9.Range("D8").Value = 0.45833333333333331
$ws9.Range("E8").Value = 2
$ws9.Range("F8").Value = "Cllient meeting"
$ws9.Range("G8").Value = "Cllient meeting"
$ws9.Range("H8").Value = "Defects tesing"

$ws9.Range("C9").Value = 0.58333333333333337
$ws9.Range("D9").Value = 0.66666666666666663
$ws9.Range("E9").Value = 2
$ws9.Range("F9").Value = "Cllient meeting"
$ws9.Range("G9").Value = "Cllient meeting"
$ws9.Range("H9").Value = "Integration testing"

$ws9.Range("C10").Value = 0.47916666666666669
$ws9.Range("D10").Value = 0.58333333333333337
$ws9.Range("E10").Value = 2.5
$ws9.Range("F10").Value = "Team Meeting"
$ws9.Range("G10").Value = "Prepare for client meeting and testing"

$ws9.Range("C11").Value = 0.54166666666666663
$ws9.Range("D11").Value = 0.41666666666666669
$ws9.Range("E11").Value = 2
$ws9.Range("F11").Value = "Cllient meeting"
$ws9.Range("G11").Value = "Cllient meeting"
$ws9.Range("H11").Value = "First draft"

$ws9.Range("H4").Select()

# ---------------------------------------------------------------------------
# Week10: clone the Week9 layout/styling, then fill its own values
# ---------------------------------------------------------------------------
$ws9.Copy($null, $ws9)
$ws10 = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws10.Name = "Week10"
$ws10.PageSetup.PrintArea = '$A$1:$H$14'
$ws10.Names.Add("Week_Start", "=Week10!`$C`$4")

$ws10.Range("H3").Value = 10

$ws10.Range("B6").Value = 44333
$ws10.Range("C6").Value = 0.41666666666666669
$ws10.Range("D6").Value = 0.625
$ws10.Range("E6").Value = 4
$ws10.Range("F6").Value = "Develop Timeline page"
$ws10.Range("G6").Value = "milestone 2 task"
$ws10.Range("H6").Value = "Frontend in progress 50%"

$ws10.Range("B7").Value = 44334
$ws10.Range("C7").Value = 0.54166666666666663
$ws10.Range("D7").Value = 0.875
$ws10.Range("E7").Value = 3
$ws10.Range("F7").Value = "Develop Timeline page"
$ws10.Range("G7").Value = "milestone 2 task"
$ws10.Range("H7").Value = "Frontend in progress 90%"

$ws10.Range("B8").Value = 44335
$ws10.Range("C8").Value = 0.375
$ws10.Range("D8").Value = 0
$ws10.Range("E8").Value = 3
$ws10.Range("F8").Value = "Work on timeline layout"
$ws10.Range("G8").Value = $null
$ws10.Range("H8").Value = $null

$ws10.Range("B9").Value = 44336
$ws10.Range("C9").Value = 0.41666666666666669
$ws10.Range("D9").Value = 0.66666666666666663
$ws10.Range("E9").Value = 6
$ws10.Range("F9").Value = "Develop Timeline page"
$ws10.Range("G9").Value = "milestone 2 task"
$ws10.Range("H9").Value = "Frontend finished 60%"

$ws10.Range("B10").Value = 44337
$ws10.Range("C10").Value = 0.79166666666666663
$ws10.Range("D10").Value = 0.875
$ws10.Range("E10").Value = 2
$ws10.Range("F10").Value = "Team Meeting"
$ws10.Range("G10").Value = "Present mailstone 1 development "
$ws10.Range("H10").Value = "Learn Angular + Ionic, finish basic routing and start hpme page"

$ws10.Range("B11").Value = 44338
$ws10.Range("C11").Value = 0.54166666666666663
$ws10.Range("D11").Value = 0.625
$ws10.Range("E11").Value = 2
$ws10.Range("F11").Value = "Cllient meeting"
$ws10.Range("G11").Value = "Cllient meeting"
$ws10.Range("H11").Value = "Finish front end part description draft"

$ws10.Range("B12").Value = 44339
$ws10.Range("C12").Value = 0.54166666666666663
$ws10.Range("D12").Value = 0.54166666666666663
$ws10.Range("E12").Value = 1
$ws10.Range("F12").Value = "Prepare Client Meeting"
$ws10.Range("G12").Value = "Cllient meeting"
$ws10.Range("H12").Value = $null

$ws10.Range("H4").Select()

# ---------------------------------------------------------------------------
# Week11: clone Week10's layout/styling, then fill its (mostly empty) values
# ---------------------------------------------------------------------------
$ws10.Copy($null, $ws10)
$ws11 = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws11.Name = "Week11"
$ws11.PageSetup.PrintArea = '$A$1:$H$14'
$ws11.Names.Add("Week_Start", "=Week11!`$C`$4")

$ws11.Range("H3").Value = 11

$ws11.Range("B6").Value = 44340
$ws11.Range("C6").Value = $null
$ws11.Range("D6").Value = $null
$ws11.Range("E6").Value = 2
$ws11.Range("F6").Value = "Cllient meeting"
$ws11.Range("G6").Value = "Cllient meeting"
$ws11.Range("H6").Value = "Finish"

$ws11.Range("B7").Value = 44341
$ws11.Range("C7").Value = $null
$ws11.Range("D7").Value = $null
$ws11.Range("E7").Value = $null
$ws11.Range("F7").Value = "Client meeting "
$ws11.Range("G7").Value = $null
$ws11.Range("H7").Value = $null

$ws11.Range("B8").Value = 44342
$ws11.Range("C8").Value = $null
$ws11.Range("D8").Value = $null
$ws11.Range("E8").Value = $null
$ws11.Range("F8").Value = $null
$ws11.Range("G8").Value = $null
$ws11.Range("H8").Value = $null

$ws11.Range("B9").Value = 44343
$ws11.Range("C9").Value = $null
$ws11.Range("D9").Value = $null
$ws11.Range("E9").Value = $null
$ws11.Range("F9").Value = $null
$ws11.Range("G9").Value = $null
$ws11.Range("H9").Value = $null

$ws11.Range("B10").Value = 44344
$ws11.Range("C10").Value = $null
$ws11.Range("D10").Value = $null
$ws11.Range("E10").Value = $null
$ws11.Range("F10").Value = $null
$ws11.Range("G10").Value = $null
$ws11.Range("H10").Value = $null

$ws11.Range("B11").Value = 44345
$ws11.Range("C11").Value = $null
$ws11.Range("D11").Value = $null
$ws11.Range("E11").Value = $null
$ws11.Range("F11").Value = $null
$ws11.Range("G11").Value = $null
$ws11.Range("H11").Value = $null

$ws11.Range("B12").Value = 44346
$ws11.Range("C12").Value = $null
$ws11.Range("D12").Value = $null
$ws11.Range("E12").Value = $null
$ws11.Range("F12").Value = $null
$ws11.Range("G12").Value = $null
$ws11.Range("H12").Value = $null

$ws11.Range("H4").Select()

# ---------------------------------------------------------------------------
# Make Week10 the active / selected tab, matching the saved workbook state
# ---------------------------------------------------------------------------
$ws10.Activate()
